$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 410.83334
$ws.Range("I28").Value = 433.4
$ws.Range("J28").Value = 298
$ws.Range("K28").Value = 433.4
$ws.Range("L28").Value = 298
$ws.Range("M28").Value = 51.60000000000002
$ws.Range("N28").Value = -1268
$ws.Range("H62").Value = 3226
$ws.Range("I62").Value = 3301.3333
$ws.Range("K62").Value = 3301.3333
$ws.Range("M62").Value = -2677.3333
$ws.Range("H65").Value = 3226
$ws.Range("I65").Value = 3301.3333
$ws.Range("K65").Value = 16506.6665
$ws.Range("M65").Value = -13386.6665
$ws.Range("H92").Value = 463.43478
$ws.Range("I92").Value = 370.45
$ws.Range("K92").Value = 370.45
$ws.Range("M92").Value = 877.55
$ws.Range("H98").Value = 627.8276
$ws.Range("I98").Value = 581.5925999999999
$ws.Range("K98").Value = 581.5925999999999
$ws.Range("M98").Value = 916.4074000000001
$ws.Range("H103").Value = 227558.64
$ws.Range("I103").Value = 278010
$ws.Range("J103").Value = 527.5
$ws.Range("K103").Value = 834030
$ws.Range("L103").Value = 1582.5
$ws.Range("M103").Value = -833444
$ws.Range("N103").Value = -2754.5
$ws.Range("H107").Value = 971.1739
$ws.Range("I107").Value = 1099.3125
$ws.Range("K107").Value = 1099.3125
$ws.Range("M107").Value = 820.6875
$ws.Range("H116").Value = 3800.6316
$ws.Range("I116").Value = 1925
$ws.Range("K116").Value = 1925
$ws.Range("M116").Value = 1517
$ws.Range("H122").Value = 627.8276
$ws.Range("I122").Value = 581.5925999999999
$ws.Range("K122").Value = 1744.7778
$ws.Range("M122").Value = 705.2222000000002
$ws.Range("H129").Value = 820.6279
$ws.Range("I129").Value = 671
$ws.Range("J129").Value = 849.7222
$ws.Range("K129").Value = 2013
$ws.Range("L129").Value = 2549.1666
$ws.Range("M129").Value = 2987
$ws.Range("N129").Value = -12549.1666
$ws.Range("H132").Value = 2173.4055
$ws.Range("I132").Value = 2240.457
$ws.Range("K132").Value = 6721.370999999999
$ws.Range("M132").Value = -4191.370999999999
$ws.Range("H137").Value = 45721.695
$ws.Range("I137").Value = 2343.4285
$ws.Range("J137").Value = 113199
$ws.Range("K137").Value = 7030.2855
$ws.Range("L137").Value = 339597
$ws.Range("M137").Value = -4480.2855
$ws.Range("N137").Value = -344697
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20064.49
$ws.Range("I32").Value = 22465.084
$ws.Range("J32").Value = 3603.2856
$ws.Range("K32").Value = 22465.084
$ws.Range("L32").Value = 3603.2856
$ws.Range("M32").Value = -22178.084
$ws.Range("N32").Value = -4177.2856
$ws.Range("H102").Value = 3552.5454
$ws.Range("J102").Value = 5416
$ws.Range("L102").Value = 5416
$ws.Range("N102").Value = -8660
$ws.Range("H132").Value = 19708.207
$ws.Range("I132").Value = 1877.7778
$ws.Range("K132").Value = 5633.3334
$ws.Range("M132").Value = -3103.3334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2392.077
$ws.Range("I99").Value = 1929.7142
$ws.Range("K99").Value = 1929.7142
$ws.Range("M99").Value = -431.7141999999999
$ws.Range("H107").Value = 1380.7693
$ws.Range("I107").Value = 881.125
$ws.Range("J107").Value = 2180.2
$ws.Range("K107").Value = 881.125
$ws.Range("L107").Value = 2180.2
$ws.Range("M107").Value = 1038.875
$ws.Range("N107").Value = -6020.2
$ws.Range("H134").Value = 71819.47
$ws.Range("I134").Value = 82522.46000000001
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 247567.38
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -245032.38
$ws.Range("N134").Value = -11820
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4139.44
$ws.Range("I99").Value = 3303.9412
$ws.Range("K99").Value = 3303.9412
$ws.Range("M99").Value = -1805.9412
$ws.Range("H126").Value = 4139.44
$ws.Range("I126").Value = 3303.9412
$ws.Range("K126").Value = 9911.8236
$ws.Range("M126").Value = -7441.8236
$ws.Range("H134").Value = 1279.909
$ws.Range("I134").Value = 1011.4286
$ws.Range("K134").Value = 3034.2858
$ws.Range("M134").Value = -499.2857999999997
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 757.5
$ws.Range("I34").Value = 515
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 1545
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1461
$ws.Range("N34").Value = -3168
$ws.Range("H131").Value = 780.46
$ws.Range("J131").Value = 781.2727
$ws.Range("L131").Value = 2343.8181
$ws.Range("N131").Value = -12423.8181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3200.182
$ws.Range("I113").Value = 2533.6667
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2533.6667
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -363.6667000000002
$ws.Range("N113").Value = -8340
$ws.Range("H132").Value = 129217.664
$ws.Range("I132").Value = 115957
$ws.Range("J132").Value = 168999.67
$ws.Range("K132").Value = 347871
$ws.Range("L132").Value = 506999.01
$ws.Range("M132").Value = -345341
$ws.Range("N132").Value = -512059.01
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 3900
$ws.Range("J26").Value = 3900
$ws.Range("L26").Value = 3900
$ws.Range("N26").Value = -4490
$ws.Range("H40").Value = 2439.4285
$ws.Range("I40").Value = 1825.1765
$ws.Range("J40").Value = 5050
$ws.Range("K40").Value = 1825.1765
$ws.Range("L40").Value = 5050
$ws.Range("M40").Value = -1689.1765
$ws.Range("N40").Value = -5322
$ws.Range("H122").Value = 1512042.1
$ws.Range("I122").Value = 3270891.2
$ws.Range("K122").Value = 9812673.600000001
$ws.Range("M122").Value = -9810223.600000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3379868.2
$ws.Range("I113").Value = 1903
$ws.Range("J113").Value = 13513764
$ws.Range("K113").Value = 5709
$ws.Range("L113").Value = 40541292
$ws.Range("M113").Value = -3539
$ws.Range("N113").Value = -40545632
$ws.Range("H122").Value = 1924.8334
$ws.Range("I122").Value = 1683.7778
$ws.Range("K122").Value = 5051.3334
$ws.Range("M122").Value = -2601.3334
$ws.Range("H136").Value = 949941.4
$ws.Range("I136").Value = 1403387.1
$ws.Range("J136").Value = 1827.5454
$ws.Range("K136").Value = 4210161.300000001
$ws.Range("L136").Value = 5482.6362
$ws.Range("M136").Value = -4207611.300000001
